$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value: (empty) -> Alvearie Team
$ws.Range("B9").Value = "Alvearie Team"

# Remove the duplicate "Contact" row (row 11), which shifts rows 12-15 up to 11-14
$ws.Range("A11").EntireRow.Delete()

# The old "Contact" / "No display for ContactDetail" row (now row 10) becomes
# "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
